$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'27.291.35"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.64%  '
$ws.Range('D3').Value = "'1.904.72"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.52%  '
$ws.Range('D4').Value = "'0.9999"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.30%  '
$ws.Range('D5').Value = "'306.42"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E6').Value = '  -0.18%  '
$ws.Range('D7').Value = "'0.5421"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +4.17%  '
$ws.Range('D8').Value = "'0.3811"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.28%  '
$ws.Range('D9').Value = "'0.07297"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.46%  '
$ws.Range('D10').Value = "'22.17"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +4.81%  '
$ws.Range('D11').Value = "'0.9026"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.29%  '
$ws.Range('D12').Value = "'0.08187"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.03%  '
$ws.Range('D13').Value = "'95.57"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.83%  '
$ws.Range('D14').Value = "'5.358"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.03%  '
$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').Value = "'1.435.81"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -24.68%  '
$ws.Range('B16').Value = 'BinanceUSD'
$ws.Range('C16').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D16').Value = "'1.001"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.22%  '
$ws.Range('B17').Value = 'Avalanche'
$ws.Range('C17').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D17').Value = "'14.89"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.00%  '
$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D18').Value = "'0.000008652"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.55%  '
$ws.Range('D19').Value = "'0.9991"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.37%  '
$ws.Range('D20').Value = "'27.302.64"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.59%  '
$ws.Range('D21').Value = "'5.053"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.54%  '
$ws.Range('D23').Value = "'6.523"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.59%  '
$ws.Range('D24').Value = "'148.63"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.23%  '
$ws.Range('D25').Value = "'2.310"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.32%  '
$ws.Range('E26').Value = '  +0.92%  '
$ws.Range('D27').Value = "'1.752"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.89%  '
$ws.Range('D28').Value = "'116.62"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.23%  '
$ws.Range('D29').Value = "'4.847"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.02%  '
$ws.Range('D30').Value = "'4.677"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.77%  '
$ws.Range('D31').Value = "'0.09187"
$ws.Range('D31').Style = 'Normal'
$ws.Range('D32').Value = "'0.8278"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +4.37%  '
$ws.Range('D33').Value = "'0.05069"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.01%  '
$ws.Range('E34').Value = '  +0.80%  '
$ws.Range('D35').Value = "'3.007"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.39%  '
$ws.Range('D36').Value = "'3.319"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.62%  '
$ws.Range('D37').Value = "'2.691"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.15%  '
$ws.Range('D38').Value = "'0.5999"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +4.96%  '
$ws.Range('D39').Value = "'0.01998"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.00%  '
$ws.Range('D40').Value = "'1.076"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.10%  '
$ws.Range('D41').Value = "'9.279"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.80%  '
$ws.Range('D42').Value = "'6.662"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.54%  '
$ws.Range('D43').Value = "'116.14"
$ws.Range('D43').Style = 'Normal'
$ws.Range('D44').Value = "'0.5140"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +5.56%  '
$ws.Range('E45').Value = '  +1.17%  '
$ws.Range('D46').Value = "'10.24"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.90%  '
$ws.Range('E47').Value = '  -0.21%  '
$ws.Range('E48').Value = '  +1.02%  '
$ws.Range('D49').Value = "'38.11"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.44%  '
$ws.Range('D50').Value = "'0.06100"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.95%  '
$ws.Range('D51').Value = "'63.39"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.40%  '
